$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("December", 16, "23:30:07", -1.69, 0, 83, 1036, 0.89),
    @("December", 17, "00:00:05", -1.81, 0, 57, 1036, 1.32),
    @("December", 17, "00:30:05", -1.88, 0, 58, 1036, 0.91),
    @("December", 17, "01:00:05", -2.13, 0, 59, 1036, 0.93),
    @("December", 17, "01:30:04", -2.31, 0, 58, 1035, 1.39),
    @("December", 17, "02:00:05", -2.69, 0, 58, 1035, 1.39),
    @("December", 17, "02:30:05", -2.88, 0, 60, 1035, 1.59),
    @("December", 17, "03:30:04", -3.19, 0, 88, 1035, 1.43),
    @("December", 17, "04:00:09", -3.19, 0, 61, 1035, 1.43),
    @("December", 17, "04:30:05", -3.13, 0, 61, 1035, 1.55),
    @("December", 17, "05:00:05", -2.88, 0, 84, 1035, 0.45),
    @("December", 17, "05:30:05", -2.81, 0, 84, 1035, 0.89),
    @("December", 17, "06:00:05", -2.88, 0, 61, 1035, 1.34),
    @("December", 17, "06:30:05", -2.63, 0, 61, 1034, 1.49),
    @("December", 17, "07:00:04", -2, 0, 58, 1034, 1.4),
    @("December", 17, "07:30:05", -2.31, 0, 58, 1034, 0.9),
    @("December", 17, "08:00:05", -2.31, 0, 58, 1034, 0.9),
    @("December", 17, "08:30:04", -2.19, 0, 56, 1034, 1.19),
    @("December", 17, "09:00:05", -2, 0, 56, 1034, 1.19),
    @("December", 17, "09:30:04", -1.38, 0, 53, 1035, 1.38),
    @("December", 17, "10:00:05", -0.31, 0, 84, 1035, 1.38),
    @("December", 17, "10:30:05", 0.6899999999999999, 0, 47, 1034, 1.03),
    @("December", 17, "11:00:04", 1.63, 0, 47, 1034, 1.03),
    @("December", 17, "11:30:05", 2.38, 0, 43, 1034, 1.09),
    @("December", 17, "12:00:05", 3.19, 0, 43, 1034, 1.09),
    @("December", 17, "12:30:04", 4.19, 0, 76, 1033, 0.89),
    @("December", 17, "13:00:05", 5.38, 0, 71, 1032, 1.48),
    @("December", 17, "13:30:05", 6.38, 0, 68, 1031, 1.49),
    @("December", 17, "14:00:04", 7.19, 0, 67, 1031, 1.79),
    @("December", 17, "14:30:05", 7.88, 0, 65, 1030, 1.79),
    @("December", 17, "15:00:05", 8.5, 0, 63, 1030, 2.24),
    @("December", 17, "15:30:04", 8.880000000000001, 0, 61, 1030, 1.79),
    @("December", 17, "16:00:05", 9.130000000000001, 0, 60, 1030, 2.24),
    @("December", 17, "16:30:04", 8.81, 0, 62, 1030, 1.34),
    @("December", 17, "17:00:05", 8.130000000000001, 0, 62, 1030, 0.89),
    @("December", 17, "17:30:05", 7.13, 0, 64, 1030, 1.34),
    @("December", 17, "18:00:05", 5.88, 0, 66, 1030, 1.75),
    @("December", 17, "18:30:05", 5, 0, 68, 1030, 0.89),
    @("December", 17, "19:00:04", 4.31, 0, 69, 1031, 1.79),
    @("December", 17, "19:30:05", 1.69, 0, 69, 1031, 2.68),
    @("December", 17, "19:40:59", 1.38, 0, 70, 1031, 1.34)
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
    $ws.Cells.Item($r, 8).Value = $rowVals[7]
}

Write-Host "Done adding rows"